# Updates the Price (D) and Volume(1h) (E) columns of the cryptos sheet
# to match the refreshed coinranking.com snapshot.
#
# Price/volume cells are stored as plain text in this workbook (inline
# strings), so each write forces a text NumberFormat before assigning the
# value (otherwise Excel auto-coerces strings like "19.11" into a float
# and mangles it to "19.109999999999999"), then restores the default
# "Normal" style so no stray per-cell formatting is introduced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "26.720.78"
Set-TextValue $ws.Range("E2") "  +0.50%  "
Set-TextValue $ws.Range("D3") "1.640.04"
Set-TextValue $ws.Range("E3") "  +0.04%  "
Set-TextValue $ws.Range("E4") "  +0.44%  "
Set-TextValue $ws.Range("D5") "217.63"
Set-TextValue $ws.Range("E5") "  +1.53%  "
Set-TextValue $ws.Range("D6") "0.503"
Set-TextValue $ws.Range("E6") "  +0.08%  "
Set-TextValue $ws.Range("E7") "  +0.32%  "
Set-TextValue $ws.Range("E8") "  +0.45%  "
Set-TextValue $ws.Range("D9") "0.0625"
Set-TextValue $ws.Range("E9") "  +0.16%  "
Set-TextValue $ws.Range("D10") "19.11"
Set-TextValue $ws.Range("E10") "  +0.24%  "
Set-TextValue $ws.Range("E11") "  +0.12%  "
Set-TextValue $ws.Range("D12") "1.868.92"
Set-TextValue $ws.Range("E12") "  -0.04%  "
Set-TextValue $ws.Range("D13") "1.633.15"
Set-TextValue $ws.Range("E13") "  -1.95%  "
Set-TextValue $ws.Range("D14") "4.15"
Set-TextValue $ws.Range("E14") "  -0.20%  "
Set-TextValue $ws.Range("E15") "  -0.19%  "
Set-TextValue $ws.Range("D16") "64.65"
Set-TextValue $ws.Range("E16") "  -0.17%  "
Set-TextValue $ws.Range("D17") "26.726.39"
Set-TextValue $ws.Range("E17") "  +0.39%  "
Set-TextValue $ws.Range("D18") "0.0₃0732"
Set-TextValue $ws.Range("E18") "  -0.99%  "
Set-TextValue $ws.Range("D19") "214.79"
Set-TextValue $ws.Range("E19") "  -0.02%  "
Set-TextValue $ws.Range("E20") "  +0.32%  "
Set-TextValue $ws.Range("E21") "  +1.08%  "
Set-TextValue $ws.Range("D22") "2.37"
Set-TextValue $ws.Range("E22") "  +7.01%  "
Set-TextValue $ws.Range("D23") "6.22"
Set-TextValue $ws.Range("E23") "  -0.08%  "
Set-TextValue $ws.Range("E24") "  -1.67%  "
Set-TextValue $ws.Range("D25") "145.31"
Set-TextValue $ws.Range("E25") "  +0.37%  "
Set-TextValue $ws.Range("E26") "  +0.40%  "
Set-TextValue $ws.Range("E27") "  -0.64%  "
Set-TextValue $ws.Range("D28") "7.17"
Set-TextValue $ws.Range("E28") "  +0.83%  "
Set-TextValue $ws.Range("D29") "15.63"
Set-TextValue $ws.Range("E29") "  -0.20%  "
Set-TextValue $ws.Range("D30") "0.0508"
Set-TextValue $ws.Range("E30") "  -0.65%  "
Set-TextValue $ws.Range("E31") "  +1.63%  "
Set-TextValue $ws.Range("E32") "  +1.01%  "
Set-TextValue $ws.Range("E33") "  +0.68%  "
Set-TextValue $ws.Range("D34") "1.285.14"
Set-TextValue $ws.Range("E34") "  +0.88%  "
Set-TextValue $ws.Range("D35") "1.53"
Set-TextValue $ws.Range("E35") "  +0.44%  "
Set-TextValue $ws.Range("E36") "  +1.35%  "
Set-TextValue $ws.Range("E37") "  +0.04%  "
Set-TextValue $ws.Range("D38") "0.537"
Set-TextValue $ws.Range("E38") "  +1.65%  "
Set-TextValue $ws.Range("D39") "0.817"
Set-TextValue $ws.Range("E39") "  -0.67%  "
Set-TextValue $ws.Range("E40") "  +0.39%  "
Set-TextValue $ws.Range("D41") "0.804"
Set-TextValue $ws.Range("E41") "  -0.42%  "
Set-TextValue $ws.Range("E42") "  -0.89%  "
Set-TextValue $ws.Range("E43") "  -2.28%  "
Set-TextValue $ws.Range("D44") "1.779.37"
Set-TextValue $ws.Range("E44") "  +0.05%  "
Set-TextValue $ws.Range("D45") "60.98"
Set-TextValue $ws.Range("E45") "  +3.11%  "
Set-TextValue $ws.Range("D46") "91.74"
Set-TextValue $ws.Range("E46") "  +0.81%  "
Set-TextValue $ws.Range("E47") "  +0.65%  "
Set-TextValue $ws.Range("E48") "  +0.30%  "
Set-TextValue $ws.Range("D49") "7.62"
Set-TextValue $ws.Range("E49") "  -1.27%  "
Set-TextValue $ws.Range("D50") "0.0965"
Set-TextValue $ws.Range("E50") "  +0.64%  "
Set-TextValue $ws.Range("E51") "  +0.06%  "
